$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right, 5 twips of space) to the
# first paragraph's pPr, matching the pBdr already present on the third
# paragraph.
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5

# Bump the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Collapse the two runs ("**ID__AFFARS_5350_topic_6__ID**" + " ") into a
# single run with the updated placeholder text (no trailing space).
$d.Content.Find.Execute("**ID__AFFARS_5350_topic_6__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5350_102_2__ID**", 2)
